$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string bookkeeping -------------------------------------------------
# The old F6 string "扩散吞噬" (shared-string index 12) must end up replaced by
# "key" (so "key" reuses slot 12, matching the target sharedStrings order),
# and F6 itself must get the *new* string "扩散吞噬1级" (which then lands on
# slot 13). Doing it in this order lets the engine's string-pool compaction
# produce exactly that layout:
#   12 -> key
#   13 -> 扩散吞噬1级
#   14 -> 扩散吞噬2级
#   15 -> 扩散吞噬3级
#   16 -> 扩散吞噬4级
$ws.Range("F6").Value = "key"

# New header cells C2/D2 ("key") - copy E2's number format/alignment (style
# index 2) across so they match the existing header look, then set the text.
$ws.Range("E2").Copy()
$ws.Range("C2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C2").Value = "key"
$ws.Range("D2").Value = "key"

# Now give F6 its real value - a brand new string, appended right after "key".
$ws.Range("F6").Value = "扩散吞噬1级"

# --- New data rows 7-9 ----------------------------------------------------
$ws.Range("C7").Value = 1001
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "测试说明"
$ws.Range("F7").Value = "扩散吞噬2级"
$ws.Range("G7").Value = "强力攻击技能"

$ws.Range("C8").Value = 1001
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = "测试说明"
$ws.Range("F8").Value = "扩散吞噬3级"
$ws.Range("G8").Value = "强力攻击技能"

$ws.Range("C9").Value = 1001
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = "测试说明"
$ws.Range("F9").Value = "扩散吞噬4级"
$ws.Range("G9").Value = "强力攻击技能"

# --- Column F width tweak (close to the recorded autofit result) ----------
$ws.Columns("F:F").ColumnWidth = 9.74

# --- Selection moves to H12, matching the saved view ------------------------
$ws.Range("H12").Select()
